$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal TEXT into a cell, even when it looks like
# a number/date (e.g. "1990-07-08") or is the empty string. Prefixing with a
# leading apostrophe forces Excel's text interpretation instead of letting it
# auto-convert to a date/number; ClearFormats() then strips the "quote
# prefix" text style that the apostrophe entry applies, leaving the cell on
# the default (General) style, matching plain data entry.
function Set-Text($cell, $text) {
    $cell.Value = "'" + $text
    $cell.ClearFormats()
}

# Row 2 becomes entirely blank (A2:I2), each cell holding an empty string.
for ($col = 1; $col -le 9; $col++) {
    Set-Text $ws.Cells.Item(2, $col) ""
}

# Row 3 gets the "P1003 / Charlie White" record (previously at row 4).
$row3 = @("P1003", "Charlie White", "charlie.white@example.com", "1990-07-08", "Male", "O-")
for ($col = 1; $col -le 6; $col++) {
    Set-Text $ws.Cells.Item(3, $col) $row3[$col - 1]
}
for ($col = 7; $col -le 9; $col++) {
    Set-Text $ws.Cells.Item(3, $col) ""
}

# Row 4 gets the "P1001 / Alice Brown" record (previously at row 2).
$row4 = @("P1001", "Alice Brown", "alice.brown@example.com", "1980-05-14", "Female", "A+")
for ($col = 1; $col -le 6; $col++) {
    Set-Text $ws.Cells.Item(4, $col) $row4[$col - 1]
}
for ($col = 7; $col -le 9; $col++) {
    Set-Text $ws.Cells.Item(4, $col) ""
}

# Row 5 (new row) gets the "P1002 / Bob Stone" record (previously at row 3).
$row5 = @("P1002", "Bob Stone", "bob.stone@example.com", "1975-11-22", "Male", "B+")
for ($col = 1; $col -le 6; $col++) {
    Set-Text $ws.Cells.Item(5, $col) $row5[$col - 1]
}
for ($col = 7; $col -le 9; $col++) {
    Set-Text $ws.Cells.Item(5, $col) ""
}

Write-Output "done"
